$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new instruction/example text + new stimulus file ---
# Leading "'" forces the quote-prefix style bit to be retained (cells B2/C2
# keep their original wrap+quotePrefix cell style) without becoming part of
# the stored text.
$ws.Range("B2").Formula = "'You will listen to sentences like this one. `n"
$ws.Range("C2").Formula = "'Sometimes the sentences will be louder, like this:"
$ws.Range("D2").Value   = "sent_stim_122220\102_SM_SNR2_0.5.wav"

# --- Row 3: new instruction/example text + new stimulus file ---
$ws.Range("B3").Formula = "'If you need to adjust the volume to make the sound louder, please do so now.`nAfter you have found a comfortable volume level, stop adjusting the volume. "
$ws.Range("C3").Formula = "'It may be difficult to understand what is being said, but please do your best.`n'Listen carefully to each sentence and pay attention to *who* is performing the action."
$ws.Range("D3").Value   = "sent_stim_122220\102_SM_SNR2.wav"

# Row 3 wraps across three visual lines with the new, longer text, so Excel
# grows the row to fit (the engine's own autofit heuristic undershoots here).
$ws.Rows(3).RowHeight = 47.25

# --- Rows 4-9: only the stimulus-file column (D) changes, to the new
# "122220" stimulus set ---
$ws.Range("D4").Value = "sent_stim_122220\103_SM_SNR2.wav"
$ws.Range("D5").Value = "sent_stim_122220\103_SF_SNR2_0.5.wav"
$ws.Range("D6").Value = "sent_stim_122220\104_OF_SNR2_0.5.wav"
$ws.Range("D7").Value = "sent_stim_122220\104_OM_SNR2.wav"
$ws.Range("D8").Value = "sent_stim_122220\105_OF_SNR2.wav"
$ws.Range("D9").Value = "sent_stim_122220\105_SF_SNR2_0.5.wav"

# --- Selection moved to B3 with the view scrolled one column right ---
$ws.Range("B3").Select()
